# Update the "through" date from 2022-09-07 to 2022-09-08 (commit: "Add data for 2022-09-16")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the sheet tab
$ws.Name = "Through 2022-09-08"

# 2. Update the header label for the current ("through") month column
$ws.Range("B1").Value = "September 2022 (through September 08)"

# 3. Apply the new carjacking counts (cumulative monthly totals bumped by the
#    new incident added on 2022-09-16)

# Row 3 - Garfield Park
$ws.Range("T3").Value = 1
$ws.Range("AC3").Value = 1
$ws.Range("AL3").Value = 4
$ws.Range("AU3").Value = 2

# Row 7 - Little Italy, UIC
$ws.Range("B7").Value = 1

# Row 8 - Humboldt Park
$ws.Range("B8").Value = 2
$ws.Range("T8").Value = 2

# Row 9 - West Town
$ws.Range("AC9").Value = 2

# Row 10 - North Lawndale
$ws.Range("T10").Value = 5

# Row 11 - West Pullman
$ws.Range("AC11").Value = 2

# Row 13 - Gage Park
$ws.Range("BM13").Value = 2

# Row 15 - Chatham
$ws.Range("AL15").Value = 1

# Row 18 - Ashburn
$ws.Range("K18").Value = 2

# Row 23 - South Chicago
$ws.Range("AU23").Value = 1

# Row 60 - Armour Square
$ws.Range("K60").Value = 1

# Row 89 - Oakland
$ws.Range("B89").Value = 2

# Row 97 - Washington Park
$ws.Range("K97").Value = 1
